$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 321
$ws.Range("J17").Value = 321
$ws.Range("L17").Value = 963
$ws.Range("N17").Value = -1299
$ws.Range("H18").Value = 2333.3333
$ws.Range("I18").Value = 2000
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -1716
$ws.Range("N18").Value = -3568
$ws.Range("H33").Value = 442
$ws.Range("J33").Value = 2300
$ws.Range("L33").Value = 2300
$ws.Range("N33").Value = -2758
$ws.Range("H40").Value = 5167
$ws.Range("J40").Value = 5000.25
$ws.Range("L40").Value = 5000.25
$ws.Range("N40").Value = -5350.25
$ws.Range("H51").Value = 10681.546
$ws.Range("I51").Value = 8928.286
$ws.Range("K51").Value = 8928.286
$ws.Range("M51").Value = -8444.286
$ws.Range("H70").Value = 3999
$ws.Range("J70").Value = 4427.4287
$ws.Range("L70").Value = 13282.2861
$ws.Range("N70").Value = -13822.2861
$ws.Range("H73").Value = 3999
$ws.Range("J73").Value = 4427.4287
$ws.Range("L73").Value = 13282.2861
$ws.Range("N73").Value = -15154.2861
$ws.Range("H125").Value = 2400.8333
$ws.Range("J125").Value = 2481
$ws.Range("L125").Value = 22329
$ws.Range("N125").Value = -27249
$ws.Range("H132").Value = 1111.5714
$ws.Range("I132").Value = 1111.5714
$ws.Range("K132").Value = 3334.7142
$ws.Range("M132").Value = -804.7142000000003
$ws.Range("H135").Value = 599.8
$ws.Range("I135").Value = 499.75
$ws.Range("K135").Value = 4497.75
$ws.Range("M135").Value = -1962.75
$ws.Range("H138").Value = 9525.454
$ws.Range("J138").Value = 7578
$ws.Range("L138").Value = 22734
$ws.Range("N138").Value = -33014
$ws.Range("H140").Value = 80708
$ws.Range("J140").Value = 80707
$ws.Range("L140").Value = 80707
$ws.Range("N140").Value = -91067
$ws.Range("H141").Value = 3049.5
$ws.Range("I141").Value = 2500
$ws.Range("J141").Value = 3599
$ws.Range("K141").Value = 7500
$ws.Range("L141").Value = 10797
$ws.Range("M141").Value = -2320
$ws.Range("N141").Value = -21157

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 6000
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = -1827
$ws.Range("N6").Value = -10346

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4959.5
$ws.Range("I86").Value = 4989.4287
$ws.Range("J86").Value = 4750
$ws.Range("K86").Value = 4989.4287
$ws.Range("L86").Value = 4750
$ws.Range("M86").Value = -3866.4287
$ws.Range("N86").Value = -6996
$ws.Range("H89").Value = 4959.5
$ws.Range("I89").Value = 4989.4287
$ws.Range("J89").Value = 4750
$ws.Range("K89").Value = 24947.1435
$ws.Range("L89").Value = 23750
$ws.Range("M89").Value = -19331.1435
$ws.Range("N89").Value = -34982

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 822
$ws.Range("I5").Value = 830.8
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 830.8
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = -718.8
$ws.Range("N5").Value = -1024
$ws.Range("H7").Value = 2726.1
$ws.Range("I7").Value = 2512.2
$ws.Range("J7").Value = 2940
$ws.Range("K7").Value = 2512.2
$ws.Range("L7").Value = 2940
$ws.Range("M7").Value = -2399.2
$ws.Range("N7").Value = -3166
$ws.Range("H25").Value = 60.5
$ws.Range("I25").Value = 60.5
$ws.Range("K25").Value = 60.5
$ws.Range("M25").Value = 113.5
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -32290
$ws.Range("H134").Value = 1904.6
$ws.Range("I134").Value = 1904.6
$ws.Range("K134").Value = 5713.799999999999
$ws.Range("M134").Value = -3178.799999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 409.81818
$ws.Range("I38").Value = 469.7143
$ws.Range("K38").Value = 1409.1429
$ws.Range("M38").Value = -1062.1429
$ws.Range("H44").Value = 2822.889
$ws.Range("I44").Value = 476.5
$ws.Range("J44").Value = 4700
$ws.Range("K44").Value = 1429.5
$ws.Range("L44").Value = 14100
$ws.Range("M44").Value = -1031.5
$ws.Range("N44").Value = -14896
$ws.Range("H131").Value = 953.19354
$ws.Range("H133").Value = 4000
$ws.Range("I133").Value = 4000
$ws.Range("K133").Value = 12000
$ws.Range("M133").Value = -6940
$ws.Range("H134").Value = 700
$ws.Range("I134").Value = 700
$ws.Range("K134").Value = 2100
$ws.Range("M134").Value = 2970
$ws.Range("H136").Value = 12498
$ws.Range("I136").Value = 12498
$ws.Range("K136").Value = 37494
$ws.Range("M136").Value = -32394
$ws.Range("H137").Value = 3991.25
$ws.Range("I137").Value = 3999.5
$ws.Range("J137").Value = 3988.5
$ws.Range("K137").Value = 11998.5
$ws.Range("L137").Value = 11965.5
$ws.Range("M137").Value = -6898.5
$ws.Range("N137").Value = -22165.5
$ws.Range("H138").Value = 3000
$ws.Range("I138").Value = 3000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9000
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -3860
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 7074.8335
$ws.Range("I139").Value = 7074.8335
$ws.Range("K139").Value = 21224.5005
$ws.Range("M139").Value = -16084.5005
$ws.Range("H140").Value = 8319.758
$ws.Range("I140").Value = 758.6667
$ws.Range("K140").Value = 2276.0001
$ws.Range("M140").Value = 2903.9999
$ws.Range("H141").Value = 9010
$ws.Range("I141").Value = 8515
$ws.Range("K141").Value = 25545
$ws.Range("M141").Value = -20365

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3430.6667
$ws.Range("I22").Value = 3896
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 3896
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -3601
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 3430.6667
$ws.Range("I27").Value = 3896
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 3896
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -3789
$ws.Range("N27").Value = -2714
$ws.Range("H43").Value = 8333.333000000001
$ws.Range("J43").Value = 8333.333000000001
$ws.Range("L43").Value = 8333.333000000001
$ws.Range("N43").Value = -8719.333000000001
$ws.Range("H46").Value = 4538.4614
$ws.Range("J46").Value = 4666.6665
$ws.Range("L46").Value = 4666.6665
$ws.Range("N46").Value = -5042.6665
$ws.Range("H55").Value = 766.25
$ws.Range("I55").Value = 549.3333
$ws.Range("K55").Value = 549.3333
$ws.Range("M55").Value = -376.3333
$ws.Range("H58").Value = 1093
$ws.Range("I58").Value = 1093
$ws.Range("K58").Value = 1093
$ws.Range("M58").Value = -833
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -4832
$ws.Range("N14").Value = -5336
$ws.Range("H51").Value = 17499.5
$ws.Range("J51").Value = 17499.5
$ws.Range("L51").Value = 17499.5
$ws.Range("N51").Value = -18519.5
$ws.Range("H86").Value = 47555
$ws.Range("J86").Value = 47555
$ws.Range("L86").Value = 47555
$ws.Range("N86").Value = -49801
$ws.Range("H89").Value = 47555
$ws.Range("J89").Value = 47555
$ws.Range("L89").Value = 237775
$ws.Range("N89").Value = -249007
$ws.Range("H100").Value = 1379.25
$ws.Range("I100").Value = 810.2143
$ws.Range("J100").Value = 5362.5
$ws.Range("K100").Value = 1620.4286
$ws.Range("L100").Value = 10725
$ws.Range("M100").Value = -1079.4286
$ws.Range("N100").Value = -11807
$ws.Range("H126").Value = 1615.6666
$ws.Range("I126").Value = 1338.8
$ws.Range("K126").Value = 4016.4
$ws.Range("M126").Value = -1546.4
$ws.Range("H127").Value = 24998
$ws.Range("J127").Value = 24998
$ws.Range("L127").Value = 24998
$ws.Range("N127").Value = -34918
$ws.Range("H132").Value = 2374.75
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 2166.3333
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -11558.9999
